# Corrige un error en la validacion de distribuciones:
# - Recalcula la fila 2 de la hoja "Data" con los nuevos numeros aleatorios.
# - Agrega las filas 3 a 11 (iteraciones 2 a 10) de la simulacion.
# - Actualiza la fila 2 (resumen de la iteracion 100) de la hoja "Ultima Iteracion".

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Hoja "Data": filas 2 a 11
# ---------------------------------------------------------------------------
$wsData = $wb.Worksheets.Item("Data")

# Extender el formato (bordes/centrado) de la fila 2 existente a las filas
# nuevas 3-11 antes de cargar los valores.
$wsData.Range("A2:M2").Copy()
$wsData.Range("A3:M11").PasteSpecial(-4122)

$dataRows = @(
    @(1,  0.9541, $false, 0.0102, "--", 0.8473, "--",    0.5372, 0, 0,   0,   0, 0),
    @(2,  0.6614, $true,  0.291,  "M",  0.274,  $false,  0.17,   0, 0,   0,   0, 0),
    @(3,  0.0648, $true,  0.7704, "M",  0.0089, $true,   0.8394, 2, 400, 400, 1, 0.3333),
    @(4,  0.4305, $true,  0.2918, "M",  0.2583, $false,  0.7566, 0, 0,   400, 1, 0.25),
    @(5,  0.0823, $true,  0.0843, "M",  0.4084, $false,  0.5346, 0, 0,   400, 1, 0.2),
    @(6,  0.3297, $true,  0.3064, "M",  0.6658, $false,  0.5375, 0, 0,   400, 1, 0.1667),
    @(7,  0.4241, $true,  0.9013, "H",  0.3023, $false,  0.7023, 0, 0,   400, 1, 0.1429),
    @(8,  0.447,  $true,  0.0415, "M",  0.3402, $false,  0.9486, 0, 0,   400, 1, 0.125),
    @(9,  0.3516, $true,  0.8421, "H",  0.6471, $false,  0.0085, 0, 0,   400, 1, 0.1111),
    @(10, 0.7778, $false, 0.2979, "--", 0.279,  "--",    0.6501, 0, 0,   400, 1, 0.1)
)

for ($i = 0; $i -lt $dataRows.Length; $i++) {
    $r = $i + 2
    $vals = $dataRows[$i]
    $wsData.Cells.Item($r, 1).Value  = $vals[0]
    $wsData.Cells.Item($r, 2).Value  = $vals[1]
    $wsData.Cells.Item($r, 3).Value  = $vals[2]
    $wsData.Cells.Item($r, 4).Value  = $vals[3]
    $wsData.Cells.Item($r, 5).Value  = $vals[4]
    $wsData.Cells.Item($r, 6).Value  = $vals[5]
    $wsData.Cells.Item($r, 7).Value  = $vals[6]
    $wsData.Cells.Item($r, 8).Value  = $vals[7]
    $wsData.Cells.Item($r, 9).Value  = $vals[8]
    $wsData.Cells.Item($r, 10).Value = $vals[9]
    $wsData.Cells.Item($r, 11).Value = $vals[10]
    $wsData.Cells.Item($r, 12).Value = $vals[11]
    $wsData.Cells.Item($r, 13).Value = $vals[12]
}

# ---------------------------------------------------------------------------
# Hoja "Ultima Iteracion": fila 2
# ---------------------------------------------------------------------------
$wsUltima = $wb.Worksheets.Item("Ultima Iteracion")

$wsUltima.Cells.Item(2, 2).Value  = 0.7398
$wsUltima.Cells.Item(2, 3).Value  = $false
$wsUltima.Cells.Item(2, 4).Value  = 0.7712
$wsUltima.Cells.Item(2, 5).Value  = "--"
$wsUltima.Cells.Item(2, 6).Value  = 0.1989
$wsUltima.Cells.Item(2, 7).Value  = "--"
$wsUltima.Cells.Item(2, 8).Value  = 0.9676
$wsUltima.Cells.Item(2, 12).Value = 8
$wsUltima.Cells.Item(2, 13).Value = 0.08

Write-Output "edit applied"
